$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.385.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.09%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.619.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.99%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'553.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.59%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'155.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.11%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  -3.05%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = "'  -4.14%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'5.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.80%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "'  -2.39%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'3.079.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.15%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'25.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.48%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'62.286.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.02%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("E16").Value = "'  -3.04%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'2.616.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -4.19%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'11.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.93%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'  -3.18%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'341.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.38%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'6.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -6.21%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "'  +0.28%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.498"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.04%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'63.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.44%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.62%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  +0.02%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'8.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.76%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'0.0₃0830"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -7.44%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = "'  +0.22%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'1.35"
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'  -3.27%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'160.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.43%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  +0.02%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("B34").Value = "'NEARProtocol"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'4.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.04%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("B35").Value = "'EthereumClassic"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'19.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.29%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  -4.86%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  -3.19%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'338.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.70%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'6.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.02%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.895"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -6.72%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = "'Filecoin"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'3.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.68%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "'OKB"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'37.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.94%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "'Mantle"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.612"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.42%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").Value = "'FirstDigitalUSD"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.12%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("B45").Value = "'InjectiveProtocol"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'20.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.67%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'2.139.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.41%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'10.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.02%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'19.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.87%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.0549"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -5.05%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'  -2.35%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "'  -3.22%  "
$ws.Range("E51").Style = "Normal"
